# Weekly data update: a new price record (dated 2023-01-17) is inserted
# right after the existing row for "Feria Lagunitas de Puerto Montt - Cebollín"
# dated 2021-05-13 (row 308), pushing all the following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 308 (shifts rows 308..386 down to 309..387)
$ws.Rows.Item(308).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A308").Value = 4
$ws.Range("B308").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C308").Value = "Los Lagos"
$ws.Range("D308").Value = 44943
$ws.Range("E308").Value = 10
$ws.Range("F308").Value = 100112037
$ws.Range("G308").Value = "Cebollín"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 160
$ws.Range("K308").Value = 6000
$ws.Range("L308").Value = 6000
$ws.Range("M308").Value = 6000
$ws.Range("N308").Value = "`$/paquete 36 unidades"
$ws.Range("O308").Value = "Región Metropolitana"
$ws.Range("P308").Value = 167
$ws.Range("Q308").Value = 36
$ws.Range("R308").Value = "Hortaliza"
